$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New vm_pu results for the 380 kV case (Case_2_91): updated slack voltage (B)
# from 1.05 pu to 1.02 pu, and recomputed bus voltages for all other buses.
$newValues = @{
    "B2" = 1.02
    "C2" = 1.036857507403616
    "D2" = 1.046780663854478
    "E2" = 1.054116006506777
    "F2" = 1.060033889811393
    "I2" = 1.040516035518783
    "J2" = 1.041963539277031
    "K2" = 1.049545140608949
    "L2" = 1.056860104761352
    "M2" = 1.062761770054698
    "N2" = 1.043443247017693
    "B3" = 1.02
    "C3" = 1.037838452201062
    "D3" = 1.047557477096394
    "E3" = 1.055050238333114
    "F3" = 1.060955706145722
    "I3" = 1.040734733412392
    "J3" = 1.042588464149037
    "K3" = 1.050133756142794
    "L3" = 1.057607206240048
    "M3" = 1.063497661035381
    "N3" = 1.044069059354697
    "B4" = 1.02
    "C4" = 1.03847339291383
    "D4" = 1.048060021635597
    "E4" = 1.055655649520737
    "F4" = 1.061552624565994
    "I4" = 1.040874576144315
    "J4" = 1.042992433698688
    "K4" = 1.050513865539294
    "L4" = 1.058090902661315
    "M4" = 1.06397364162705
    "N4" = 1.044473602587445
    "B5" = 1.02
    "C5" = 1.038740370155439
    "D5" = 1.048271264894498
    "E5" = 1.055910378798995
    "F5" = 1.061803673575277
    "I5" = 1.040932965770764
    "J5" = 1.043162166294804
    "K5" = 1.05067347965812
    "L5" = 1.058294312682875
    "M5" = 1.06417369731897
    "N5" = 1.044643576223315
    "B6" = 1.02
    "C6" = 1.03878519961518
    "D6" = 1.048306731980527
    "E6" = 1.05595316150357
    "F6" = 1.061845831901415
    "I6" = 1.040942746157349
    "J6" = 1.043190659480776
    "K6" = 1.050700268761966
    "L6" = 1.058328469846404
    "M6" = 1.064207284789912
    "N6" = 1.04467210987288
    "B7" = 1.02
    "C7" = 1.038476960088331
    "D7" = 1.048062844383858
    "E7" = 1.055659052384852
    "F7" = 1.061555978686043
    "I7" = 1.040875357922851
    "J7" = 1.042994702052926
    "K7" = 1.050515999034611
    "L7" = 1.058093620385612
    "M7" = 1.063976314965219
    "N7" = 1.044475874163006
    "B8" = 1.02
    "C8" = 1.03718897970506
    "D8" = 1.047043212671406
    "E8" = 1.054431547292069
    "F8" = 1.060345329787907
    "I8" = 1.04059029080638
    "J8" = 1.042174817640416
    "K8" = 1.049744223669719
    "L8" = 1.05711253430876
    "M8" = 1.063010506588254
    "N8" = 1.043654825420587
    "B9" = 1.02
    "C9" = 1.03492098243166
    "D9" = 1.045245734991427
    "E9" = 1.052275478659146
    "F9" = 1.058215447420282
    "I9" = 1.040075206248577
    "J9" = 1.040727057215955
    "K9" = 1.048378444068919
    "L9" = 1.055385862065512
    "M9" = 1.06130721880536
    "N9" = 1.042205009010264
    "B10" = 1.02
    "C10" = 1.033410087386366
    "D10" = 1.04404697699615
    "E10" = 1.050842844948392
    "F10" = 1.056797906346614
    "I10" = 1.039723271634587
    "J10" = 1.039759897596977
    "K10" = 1.047464068911248
    "L10" = 1.054236236732658
    "M10" = 1.060170805906565
    "N10" = 1.041236475913692
    "B11" = 1.02
    "C11" = 1.03275612145195
    "D11" = 1.043527811382191
    "E11" = 1.05022363887991
    "F11" = 1.056184674435993
    "I11" = 1.03956886021862
    "J11" = 1.039340644583111
    "K11" = 1.047067230603705
    "L11" = 1.053738801750469
    "M11" = 1.05967852829511
    "N11" = 1.040816627512457
    "B12" = 1.02
    "C12" = 1.032513249135847
    "D12" = 1.043334956891971
    "E12" = 1.049993809651554
    "F12" = 1.055956979954021
    "I12" = 1.03951120175398
    "J12" = 1.039184845911267
    "K12" = 1.046919691597312
    "L12" = 1.05355408733167
    "M12" = 1.059495644857996
    "N12" = 1.040660607588629
    "B13" = 1.02
    "C13" = 1.032565344266123
    "D13" = 1.043376325417755
    "E13" = 1.050043101023328
    "F13" = 1.056005817231681
    "I13" = 1.039523583402272
    "J13" = 1.039218268386865
    "K13" = 1.046951345337902
    "L13" = 1.053593706693885
    "M13" = 1.059534875301175
    "N13" = 1.040694077527975
    "B14" = 1.02
    "C14" = 1.032736044737052
    "D14" = 1.04351187023424
    "E14" = 1.050204637619573
    "F14" = 1.056165851354774
    "I14" = 1.039564100334015
    "J14" = 1.039327767628978
    "K14" = 1.047055037744776
    "L14" = 1.053723532069905
    "M14" = 1.059663411691861
    "N14" = 1.040803732271573
    "B15" = 1.02
    "C15" = 1.032841224211698
    "D15" = 1.0435953821437
    "E15" = 1.050304188404091
    "F15" = 1.056264465241622
    "I15" = 1.039589023992363
    "J15" = 1.039395224544654
    "K15" = 1.047118908115899
    "L15" = 1.053803529095821
    "M15" = 1.059742603308215
    "N15" = 1.040871284983806
    "B16" = 1.02
    "C16" = 1.033453494501337
    "D16" = 1.04408143038861
    "E16" = 1.050883963600487
    "F16" = 1.056838616690245
    "I16" = 1.039733476866077
    "J16" = 1.039787712265592
    "K16" = 1.047490386704309
    "L16" = 1.054269257527076
    "M16" = 1.060203472546277
    "N16" = 1.041264330082326
    "B17" = 1.02
    "C17" = 1.033837625829864
    "D17" = 1.044386290810365
    "E17" = 1.051247945613132
    "F17" = 1.057198920843772
    "I17" = 1.039823547506007
    "J17" = 1.040033785005863
    "K17" = 1.047723162923429
    "L17" = 1.054561493813851
    "M17" = 1.060492509867132
    "N17" = 1.041510752274113
    "B18" = 1.02
    "C18" = 1.03406170834437
    "D18" = 1.044564101450759
    "E18" = 1.05146035948037
    "F18" = 1.057409135346538
    "I18" = 1.039875889116401
    "J18" = 1.040177270016233
    "K18" = 1.047858849707751
    "L18" = 1.054731985003545
    "M18" = 1.060661080765654
    "N18" = 1.041654441049656
    "B19" = 1.02
    "C19" = 1.034138118953345
    "D19" = 1.044624728676106
    "E19" = 1.051532805627314
    "F19" = 1.057480822386289
    "I19" = 1.039893703150304
    "J19" = 1.040226187034627
    "K19" = 1.047905100504026
    "L19" = 1.054790123951782
    "M19" = 1.060718555719161
    "N19" = 1.041703427535827
    "B20" = 1.02
    "C20" = 1.033796409562823
    "D20" = 1.044353583124884
    "E20" = 1.051208882472396
    "F20" = 1.057160257901787
    "I20" = 1.039813903952092
    "J20" = 1.040007388385195
    "K20" = 1.047698197297882
    "L20" = 1.054530136029145
    "M20" = 1.06046150092688
    "N20" = 1.041484318167217
    "B21" = 1.02
    "C21" = 1.032685776602063
    "D21" = 1.043471956021678
    "E21" = 1.050157064375477
    "F21" = 1.056118722879514
    "I21" = 1.039552177472955
    "J21" = 1.039295524742885
    "K21" = 1.047024506651937
    "L21" = 1.05368530020441
    "M21" = 1.05962556173848
    "N21" = 1.040771443596883
    "B22" = 1.02
    "C22" = 1.031987707781839
    "D22" = 1.042917565853151
    "E22" = 1.049496737258617
    "F22" = 1.055464373264097
    "I22" = 1.03938586573419
    "J22" = 1.038847546408606
    "K22" = 1.046600146518273
    "L22" = 1.053154437540929
    "M22" = 1.05909980271817
    "N22" = 1.040322829081985
    "B23" = 1.02
    "C23" = 1.032357745300593
    "D23" = 1.043211465435525
    "E23" = 1.049846694684219
    "F23" = 1.055811208158325
    "I23" = 1.039474196816682
    "J23" = 1.039085066030751
    "K23" = 1.046825181909761
    "L23" = 1.053435827317849
    "M23" = 1.059378533408193
    "N23" = 1.040560686009238
    "B24" = 1.02
    "C24" = 1.033815033341818
    "D24" = 1.044368362348825
    "E24" = 1.051226533086552
    "F24" = 1.057177727850594
    "I24" = 1.039818262061753
    "J24" = 1.040019316021071
    "K24" = 1.047709478460276
    "L24" = 1.054544305155091
    "M24" = 1.060475512593072
    "N24" = 1.041496262741704
    "B25" = 1.02
    "C25" = 1.035507122216073
    "D25" = 1.045710508560443
    "E25" = 1.052832042876469
    "F25" = 1.058765658638641
    "I25" = 1.040209876539558
    "J25" = 1.041101691249442
    "K25" = 1.048732214227779
    "L25" = 1.055831990048814
    "M25" = 1.061747720393922
    "N25" = 1.04258017506705
}

foreach ($addr in $newValues.Keys) {
    $ws.Range($addr).Value = $newValues[$addr]
}

